$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title (H1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Book of Dead for Free - Slot
#    Game Review" right before the final (italic) paragraph.
#    To avoid inheriting the italic / list-bullet formatting of the
#    neighbouring paragraphs, we build the new paragraph in a "clean"
#    (Normal-style, non-italic) spot of the document, format it, then
#    move (cut/paste) it into its final position.
# ------------------------------------------------------------------
$plainPara = $d.Paragraphs.Item(4)
$insertAnchor = $d.Range($plainPara.Range.End - 1, $plainPara.Range.End - 1)
$insertAnchor.InsertAfter("`rPlay Book of Dead for Free - Slot Game Review")

$newPara = $d.Paragraphs.Item(5)
$newParaTextOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaTextOnly.Bold = 1

$newParaWhole = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newParaWhole.Cut()

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$targetAnchor = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$targetAnchor.Paste()

# ------------------------------------------------------------------
# 3) Replace the text of the (now last) italic paragraph with the new
#    description copy, keeping its italic formatting intact.
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$descPara = $d.Paragraphs.Item($paraCount)
$descTextOnly = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descTextOnly.Text = "Read our review of the Book of Dead slot game. Play for free with a wide range of betting options and unique features like wild and scatter symbols represented by the same symbol."
